# Auto-generated script to update LeveProfits market data values
# per scheduled runner commit (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 329.8
$ws.Range("I4").Value = 287.25
$ws.Range("K4").Value = 287.25
$ws.Range("M4").Value = -173.25
$ws.Range("H33").Value = 556.9545000000001
$ws.Range("I33").Value = 196.85715
$ws.Range("K33").Value = 196.85715
$ws.Range("M33").Value = 32.14285000000001
$ws.Range("H82").Value = 2937.25
$ws.Range("I82").Value = 583
$ws.Range("K82").Value = 1749
$ws.Range("M82").Value = -1343
$ws.Range("H85").Value = 2937.25
$ws.Range("I85").Value = 583
$ws.Range("K85").Value = 1749
$ws.Range("M85").Value = -345
$ws.Range("H86").Value = 6312
$ws.Range("J86").Value = 6312
$ws.Range("L86").Value = 6312
$ws.Range("N86").Value = -8558
$ws.Range("H88").Value = 2654.8333
$ws.Range("J88").Value = 2654.8333
$ws.Range("L88").Value = 2654.8333
$ws.Range("N88").Value = -3466.8333
$ws.Range("H89").Value = 6312
$ws.Range("J89").Value = 6312
$ws.Range("L89").Value = 31560
$ws.Range("N89").Value = -42792
$ws.Range("H91").Value = 2654.8333
$ws.Range("J91").Value = 2654.8333
$ws.Range("L91").Value = 2654.8333
$ws.Range("N91").Value = -5462.8333
$ws.Range("H97").Value = 2579.8
$ws.Range("J97").Value = 2579.8
$ws.Range("L97").Value = 7739.400000000001
$ws.Range("N97").Value = -8731.400000000001
$ws.Range("H121").Value = 999
$ws.Range("J121").Value = 999
$ws.Range("L121").Value = 2997
$ws.Range("N121").Value = -6491

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2385.5
$ws.Range("J3").Value = 6100
$ws.Range("L3").Value = 6100
$ws.Range("N3").Value = -6330
$ws.Range("H6").Value = 203795.95
$ws.Range("I6").Value = 666821
$ws.Range("J6").Value = 5356.643
$ws.Range("K6").Value = 666821
$ws.Range("L6").Value = 5356.643
$ws.Range("M6").Value = -666648
$ws.Range("N6").Value = -5702.643
$ws.Range("H8").Value = 3343816.8
$ws.Range("J8").Value = 4298.6665
$ws.Range("L8").Value = 4298.6665
$ws.Range("N8").Value = -4586.6665
$ws.Range("H74").Value = 1608
$ws.Range("I74").Value = 1571.091
$ws.Range("J74").Value = 2014
$ws.Range("K74").Value = 1571.091
$ws.Range("L74").Value = 2014
$ws.Range("M74").Value = -697.0909999999999
$ws.Range("N74").Value = -3762
$ws.Range("H77").Value = 1608
$ws.Range("I77").Value = 1571.091
$ws.Range("J77").Value = 2014
$ws.Range("K77").Value = 7855.455
$ws.Range("L77").Value = 10070
$ws.Range("M77").Value = -3487.455
$ws.Range("N77").Value = -18806

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I7").Value = 3800513.5
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 3800513.5
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = -3800400.5
$ws.Range("N7").Value = -426
$ws.Range("H10").Value = 997.5
$ws.Range("J10").Value = 997.5
$ws.Range("L10").Value = 997.5
$ws.Range("N10").Value = -1277.5
$ws.Range("H12").Value = 960.4286
$ws.Range("J12").Value = 1137.25
$ws.Range("L12").Value = 1137.25
$ws.Range("N12").Value = -1473.25
$ws.Range("H99").Value = 2393.0833
$ws.Range("I99").Value = 2252.75
$ws.Range("J99").Value = 2673.75
$ws.Range("K99").Value = 2252.75
$ws.Range("L99").Value = 2673.75
$ws.Range("M99").Value = -754.75
$ws.Range("N99").Value = -5669.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 4149.75
$ws.Range("I3").Value = 799.75
$ws.Range("J3").Value = 7499.75
$ws.Range("K3").Value = 799.75
$ws.Range("L3").Value = 7499.75
$ws.Range("M3").Value = -686.75
$ws.Range("N3").Value = -7725.75
$ws.Range("H4").Value = 6385.2
$ws.Range("J4").Value = 7624.75
$ws.Range("L4").Value = 7624.75
$ws.Range("N4").Value = -7848.75
$ws.Range("H12").Value = 29519.7
$ws.Range("I12").Value = 200
$ws.Range("K12").Value = 200
$ws.Range("M12").Value = -30

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1499.5
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7622
$ws.Range("H71").Value = 1499.5
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 18000
$ws.Range("N71").Value = -26112
$ws.Range("H139").Value = 2042.2
$ws.Range("I139").Value = 1303.25
$ws.Range("K139").Value = 3909.75
$ws.Range("M139").Value = 1230.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1419.2
$ws.Range("I80").Value = 1525
$ws.Range("J80").Value = 996
$ws.Range("K80").Value = 1525
$ws.Range("L80").Value = 996
$ws.Range("M80").Value = -527
$ws.Range("N80").Value = -2992
$ws.Range("H83").Value = 1419.2
$ws.Range("I83").Value = 1525
$ws.Range("J83").Value = 996
$ws.Range("K83").Value = 7625
$ws.Range("L83").Value = 4980
$ws.Range("M83").Value = -2633
$ws.Range("N83").Value = -14964
$ws.Range("H132").Value = 2738.5
$ws.Range("I132").Value = 2604
$ws.Range("K132").Value = 7812
$ws.Range("M132").Value = -5282

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2175.1428
$ws.Range("I22").Value = 1350.3334
$ws.Range("K22").Value = 1350.3334
$ws.Range("M22").Value = -1055.3334
$ws.Range("H27").Value = 2175.1428
$ws.Range("I27").Value = 1350.3334
$ws.Range("K27").Value = 1350.3334
$ws.Range("M27").Value = -1243.3334
$ws.Range("H46").Value = 8750
$ws.Range("J46").Value = 6666.6665
$ws.Range("L46").Value = 6666.6665
$ws.Range("N46").Value = -7042.6665
$ws.Range("H61").Value = 4998
$ws.Range("I61").Value = 4998
$ws.Range("K61").Value = 4998
$ws.Range("M61").Value = -4796
$ws.Range("H113").Value = 4998
$ws.Range("I113").Value = 4998
$ws.Range("K113").Value = 4998
$ws.Range("M113").Value = -2828
$ws.Range("H132").Value = 4117.75
$ws.Range("I132").Value = 3851.3
$ws.Range("K132").Value = 11553.9
$ws.Range("M132").Value = -9023.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 320
$ws.Range("I23").Value = 255
$ws.Range("K23").Value = 255
$ws.Range("M23").Value = -26
$ws.Range("H136").Value = 4052.182
$ws.Range("I136").Value = 3396.4285
$ws.Range("K136").Value = 10189.2855
$ws.Range("M136").Value = -7639.2855
